{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst replacements = [\n  \"\u26a1\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 03.06.24: \u26a1\ud83d\ude80\",\n  \"Better & Faster Large Language Models via Multi-token Prediction\",\n  \"\u05d0\u05ea\u05dd \u05d1\u05d8\u05d7 \u05e9\u05d9\u05d5\u05d3\u05e2\u05d9\u05dd \u05d0\u05e0\u05d5 \u05e8\u05d2\u05d9\u05dc\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d7\u05d9\u05d6\u05d5\u05d9 \u05d8\u05d5\u05e7\u05df \u05d4\u05d1\u05d0 \u05d1\u05d4\u05d9\u05e0\u05ea\u05dd \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05e7\u05d5\u05d3\u05de\u05d9\u05dd (\u05d4\u05e7\u05e9\u05e8 \u05d0\u05d5 \u05e7\u05d5\u05e0\u05d8\u05e7\u05e1\u05d8). \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 (\u05e9\u05e7\u05d9\u05d1\u05dc \u05d3\u05d9 \u05d4\u05e8\u05d1\u05d4 pr \u05db\u05e9\u05d9\u05e6\u05d0) \u05de\u05e6\u05d9\u05e2 \u05dc\u05d7\u05d6\u05d5\u05ea \u05db\u05de\u05d4 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e2\u05d5\u05e7\u05d1\u05d9\u05dd \u05d1\u05d5 \u05d6\u05de\u05e0\u05d9\u05ea \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05e7\u05e9\u05e8. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05e8\u05d0\u05d5 \u05e9\u05d6\u05d4 \u05d9\u05db\u05d5\u05dc \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc - \u05d6\u05d4 \u05dc\u05d0 \u05de\u05e4\u05ea\u05d9\u05e2(\u05dc\u05ea\u05d7\u05d5\u05e9\u05ea\u05d9) \u05db\u05d9 \u05de\u05e9\u05d9\u05de\u05ea \u05d7\u05d9\u05d6\u05d5\u05d9 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05de\u05e8\u05d5\u05d1\u05d9\u05dd \u05d3\u05d5\u05e8\u05e9\u05ea \u05de\u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d1\u05e0\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de\u05e2\u05de\u05d9\u05e7\u05d4 \u05e9\u05dc \u05d4\u05e9\u05e4\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d2\u05dd \u05e2\u05e9\u05d5\u05d9\u05d4 \u05dc\u05ea\u05e8\u05d5\u05dd \u05dc\u05d4\u05d0\u05e6\u05ea \u05d6\u05de\u05df \u05e8\u05d9\u05e6\u05d4 \u05d5\u05d4\u05e8\u05d5\u05d5\u05d7\u05d9\u05dd \u05d2\u05d3\u05dc\u05d9\u05dd \u05e2\u05dd \u05d2\u05d5\u05d3\u05dc \u05d4\u05de\u05d5\u05d3\u05dc.\",\n  \"\u05de\u05d0\u05de\u05e8: https://arxiv.org/pdf/2404.19737\",\n  \"\u05d8\u05dc\u05d2\u05e8\u05dd: https://t.me/MathyAIwithMike/69\",\n];\n\nfor (let i = 0; i < replacements.length && i < paragraphs.items.length; i++) {\n  paragraphs.items[i].insertText(replacements[i], \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$d.Paragraphs.Item(1).Range.Text = \"\u26a1\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 03.06.24: \u26a1\ud83d\ude80\"\n$d.Paragraphs.Item(2).Range.Text = \"Better & Faster Large Language Models via Multi-token Prediction\"\n$d.Paragraphs.Item(3).Range.Text = \"\u05d0\u05ea\u05dd \u05d1\u05d8\u05d7 \u05e9\u05d9\u05d5\u05d3\u05e2\u05d9\u05dd \u05d0\u05e0\u05d5 \u05e8\u05d2\u05d9\u05dc\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d7\u05d9\u05d6\u05d5\u05d9 \u05d8\u05d5\u05e7\u05df \u05d4\u05d1\u05d0 \u05d1\u05d4\u05d9\u05e0\u05ea\u05dd \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05e7\u05d5\u05d3\u05de\u05d9\u05dd (\u05d4\u05e7\u05e9\u05e8 \u05d0\u05d5 \u05e7\u05d5\u05e0\u05d8\u05e7\u05e1\u05d8). \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 (\u05e9\u05e7\u05d9\u05d1\u05dc \u05d3\u05d9 \u05d4\u05e8\u05d1\u05d4 pr \u05db\u05e9\u05d9\u05e6\u05d0) \u05de\u05e6\u05d9\u05e2 \u05dc\u05d7\u05d6\u05d5\u05ea \u05db\u05de\u05d4 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e2\u05d5\u05e7\u05d1\u05d9\u05dd \u05d1\u05d5 \u05d6\u05de\u05e0\u05d9\u05ea \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05e7\u05e9\u05e8. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05e8\u05d0\u05d5 \u05e9\u05d6\u05d4 \u05d9\u05db\u05d5\u05dc \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc - \u05d6\u05d4 \u05dc\u05d0 \u05de\u05e4\u05ea\u05d9\u05e2(\u05dc\u05ea\u05d7\u05d5\u05e9\u05ea\u05d9) \u05db\u05d9 \u05de\u05e9\u05d9\u05de\u05ea \u05d7\u05d9\u05d6\u05d5\u05d9 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05de\u05e8\u05d5\u05d1\u05d9\u05dd \u05d3\u05d5\u05e8\u05e9\u05ea \u05de\u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d1\u05e0\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de\u05e2\u05de\u05d9\u05e7\u05d4 \u05e9\u05dc \u05d4\u05e9\u05e4\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d2\u05dd \u05e2\u05e9\u05d5\u05d9\u05d4 \u05dc\u05ea\u05e8\u05d5\u05dd \u05dc\u05d4\u05d0\u05e6\u05ea \u05d6\u05de\u05df \u05e8\u05d9\u05e6\u05d4 \u05d5\u05d4\u05e8\u05d5\u05d5\u05d7\u05d9\u05dd \u05d2\u05d3\u05dc\u05d9\u05dd \u05e2\u05dd \u05d2\u05d5\u05d3\u05dc \u05d4\u05de\u05d5\u05d3\u05dc.\"\n$d.Paragraphs.Item(4).Range.Text = \"\u05de\u05d0\u05de\u05e8: https://arxiv.org/pdf/2404.19737\"\n$d.Paragraphs.Item(5).Range.Text = \"\u05d8\u05dc\u05d2\u05e8\u05dd: https://t.me/MathyAIwithMike/69\"\n"}
